$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F24").Value = 9
$ws.Range("G24").Value = 244.35
$ws.Range("F26").Value = 86
$ws.Range("G26").Value = 3964.6
$ws.Range("B46").Value = 28554.81
$ws.Range("F48").Value = 36
$ws.Range("G48").Value = 7083.72
$ws.Range("F52").Value = 95
$ws.Range("G52").Value = 3458.95
$ws.Range("F55").Value = 216
$ws.Range("G55").Value = 41664.24
$ws.Range("F62").Value = 9
$ws.Range("G62").Value = 504.9
$ws.Range("F77").Value = 51
$ws.Range("G77").Value = 3033.99
$ws.Range("B85").Value = 161066.38
$ws.Range("F95").Value = 12
$ws.Range("G95").Value = 3015
$ws.Range("F99").Value = 23
$ws.Range("G99").Value = 2592.56
$ws.Range("F100").Value = 9
$ws.Range("G100").Value = 1071.18
$ws.Range("B103").Value = 31654.72
$ws.Range("F120").Value = 312
$ws.Range("G120").Value = 25325.04
$ws.Range("F130").Value = 112
$ws.Range("G130").Value = 8845.76
$ws.Range("B134").Value = 84878.16
$ws.Range("F152").Value = 62
$ws.Range("G152").Value = 2630.66
$ws.Range("F157").Value = 30
$ws.Range("G157").Value = 964.2
$ws.Range("B159").Value = 74196.28999999999
$ws.Range("F175").Value = 184
$ws.Range("G175").Value = 5617.52
$ws.Range("F178").Value = 142
$ws.Range("G178").Value = 4546.84
$ws.Range("B180").Value = 37686.58
$ws.Range("F184").Value = 75
$ws.Range("G184").Value = 9003
$ws.Range("F191").Value = 62
$ws.Range("G191").Value = 7844.86
$ws.Range("B198").Value = 46549.63
$ws.Range("F218").Value = 25
$ws.Range("G218").Value = 1956.5
$ws.Range("B228").Value = 35829.07
$ws.Range("F265").Value = 42
$ws.Range("G265").Value = 2197.44
$ws.Range("B267").Value = 21830.47
$ws.Range("F314").Value = 88
$ws.Range("G314").Value = 1778.48
$ws.Range("F332").Value = 74
$ws.Range("G332").Value = 2464.2
$ws.Range("F337").Value = 28
$ws.Range("G337").Value = 7619.36
$ws.Range("B349").Value = 143429.93
$ws.Range("F355").Value = 60
$ws.Range("G355").Value = 10367.4
$ws.Range("F356").Value = 3
$ws.Range("G356").Value = 119.16
$ws.Range("F361").Value = 60
$ws.Range("G361").Value = 4422.6
$ws.Range("F367").Value = 18
$ws.Range("G367").Value = 2384.64
$ws.Range("F399").Value = 155
$ws.Range("G399").Value = 9098.5
$ws.Range("F402").Value = 22
$ws.Range("G402").Value = 1198.56
$ws.Range("F409").Value = 198
$ws.Range("G409").Value = 33923.34
$ws.Range("F418").Value = 9
$ws.Range("G418").Value = 535.23
$ws.Range("F421").Value = 345
$ws.Range("G421").Value = 13713.75
$ws.Range("B423").Value = 164111.41
$ws.Range("F426").Value = 27
$ws.Range("G426").Value = 4956.66
$ws.Range("F429").Value = 9
$ws.Range("G429").Value = 4856.49
$ws.Range("B437").Value = 24033.49
$ws.Range("F469").Value = 34
$ws.Range("G469").Value = 3161.66
$ws.Range("F479").Value = 23
$ws.Range("G479").Value = 3326.72
$ws.Range("F480").Value = 68
$ws.Range("G480").Value = 8493.879999999999
$ws.Range("B481").Value = 45251.08
$ws.Range("F483").Value = 0
$ws.Range("G483").Value = 0
$ws.Range("F496").Value = 102
$ws.Range("G496").Value = 15327.54
$ws.Range("B497").Value = 40388.43
$ws.Range("F518").Value = 70
$ws.Range("G518").Value = 7387.8
$ws.Range("F521").Value = 174
$ws.Range("G521").Value = 4753.68
$ws.Range("F526").Value = 796
$ws.Range("G526").Value = 76893.60000000001
$ws.Range("B532").Value = 156334.54
$ws.Range("F549").Value = 35
$ws.Range("G549").Value = 436.45
$ws.Range("B556").Value = 15142.14
$ws.Range("F578").Value = 0
$ws.Range("G578").Value = 0
$ws.Range("B583").Value = 41838.24
$ws.Range("F592").Value = 110
$ws.Range("G592").Value = 3443
$ws.Range("F601").Value = 43
$ws.Range("G601").Value = 1156.27
$ws.Range("F609").Value = 37
$ws.Range("G609").Value = 3634.51
$ws.Range("B610").Value = 59066.86
$ws.Range("F621").Value = 245
$ws.Range("G621").Value = 14876.4
$ws.Range("F625").Value = 35
$ws.Range("G625").Value = 2251.2
$ws.Range("B638").Value = 151852.44
$ws.Range("F647").Value = 0
$ws.Range("G647").Value = 0
$ws.Range("B650").Value = 9626.68
$ws.Range("F661").Value = 16
$ws.Range("G661").Value = 2770.72
$ws.Range("B667").Value = 27328.1
$ws.Range("F669").Value = 48
$ws.Range("G669").Value = 3803.04
$ws.Range("F672").Value = 74
$ws.Range("G672").Value = 19690.66
$ws.Range("B688").Value = 90366.61
$ws.Range("F714").Value = 74
$ws.Range("G714").Value = 13173.48
$ws.Range("F716").Value = 31
$ws.Range("G716").Value = 2943.45
$ws.Range("F717").Value = 50
$ws.Range("G717").Value = 1360
$ws.Range("F718").Value = 144
$ws.Range("G718").Value = 3916.8
$ws.Range("F719").Value = 120
$ws.Range("G719").Value = 3264
$ws.Range("B720").Value = 36015.58
$ws.Range("F745").Value = 95
$ws.Range("G745").Value = 5880.5
$ws.Range("F760").Value = 33
$ws.Range("G760").Value = 22470.36
$ws.Range("B773").Value = 149390.72
$ws.Range("F799").Value = 5
$ws.Range("G799").Value = 124.95
$ws.Range("B805").Value = 4350.11
$ws.Range("F817").Value = 3
$ws.Range("G817").Value = 244.68
$ws.Range("F820").Value = 58
$ws.Range("G820").Value = 4730.48
$ws.Range("F821").Value = 137
$ws.Range("G821").Value = 18234.7
$ws.Range("F822").Value = 47
$ws.Range("G822").Value = 5199.61
$ws.Range("F825").Value = 25
$ws.Range("G825").Value = 932
$ws.Range("B837").Value = 201080.82
$ws.Range("F839").Value = 33
$ws.Range("G839").Value = 8232.84
$ws.Range("F840").Value = 47
$ws.Range("G840").Value = 8120.19
$ws.Range("F843").Value = 80
$ws.Range("G843").Value = 8704.799999999999
$ws.Range("F844").Value = 16
$ws.Range("G844").Value = 405.28
$ws.Range("F849").Value = 44
$ws.Range("G849").Value = 6438.52
$ws.Range("F852").Value = 29
$ws.Range("G852").Value = 4140.33
$ws.Range("F859").Value = 344
$ws.Range("G859").Value = 35401.04
$ws.Range("F861").Value = 312
$ws.Range("G861").Value = 11490.96
$ws.Range("F862").Value = 21
$ws.Range("G862").Value = 991.41
$ws.Range("F863").Value = 145
$ws.Range("G863").Value = 4612.45
$ws.Range("F865").Value = 105
$ws.Range("G865").Value = 5238.45
$ws.Range("B867").Value = 212391.21
$ws.Range("F888").Value = 4
$ws.Range("G888").Value = 2207.92
$ws.Range("F897").Value = 8
$ws.Range("G897").Value = 10249.04
$ws.Range("B904").Value = 43595.52
$ws.Range("B923").Value = 2792314.51
$ws.Range("B924").Value = 2792314.51
